# The workbook's "Sheet4" tab holds the three seed/test-account names that
# the rest of the sheets (Sheet2, Sheet3, Sheet5-Sheet8) reference via
# formulas such as "=Sheet4!A2". Bumping the numeric suffix from 48 to 50
# (a new Jenkins test run id) here is enough - Excel recalculates every
# dependent formula/cached value across the workbook automatically.
$wb  = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("Sheet4")

$ws4.Range("A2").Value = "tavalinetont50"
$ws4.Range("C2").Value = "puhtaloom50"
$ws4.Range("E2").Value = "filmweird50"
